# "version con revision de estilo"
#
# 1. Update the "APROVECHADO" (sí/no) column on the "CUADERNO DEL PROFESOR"
#    sheet for several rows.
# 2. Move the active/selected tab from "GUION" to "CUADERNO DEL PROFESOR"
#    and change the selected range on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUADERNO DEL PROFESOR")

# --- Column C ("sí"/"no") updates -----------------------------------------
# Rows that flip from "no" to "sí"
$rowsToSi = @(2, 3, 4, 6, 7, 8, 9, 10, 11, 12, 13, 16, 17, 19, 20, 21)
foreach ($r in $rowsToSi) {
    $ws.Range("C$r").Value = "sí"
}

# Row that flips from "sí" to "no"
$ws.Range("C14").Value = "no"

# --- Selection / active sheet ----------------------------------------------
# Make "CUADERNO DEL PROFESOR" the active sheet and select B14:C14.
$ws.Activate()
$ws.Range("B14:C14").Select()
